# "Arreglos en el home"
# Applies 4 changes to Listado_Faltantes.docx:
#  1) paragraph "... ver donde se trunca el numero." -> wrap "numero" with gramStart/gramEnd
#     (in addition to the existing spellStart/spellEnd pair).
#  2) paragraph "... ver tema del token." -> split off "token" with spellStart/spellEnd,
#     and add a brand-new list paragraph "FILTRO POR ZONA A LA HORA DE BUSCAR." right
#     after it (the _GoBack bookmark moves onto the new paragraph).
#  3) paragraph "... Falta mostrar mas notificaciones." -> wrap "mas" with gramStart/gramEnd
#     (in addition to the existing spellStart/spellEnd pair).
#  4) heading paragraph "listado-solicitudes-ofrecidas" -> wrap the whole run with
#     gramStart/gramEnd.
#
# Strategy: locate each target paragraph, then replace its *entire* Range in one
# InsertXML call with the fully reconstructed OOXML for that paragraph (and, where a
# paragraph is being added, the following paragraph too). Replacing the whole paragraph
# range in a single shot keeps the ordering of the (zero-width) w:proofErr markers
# deterministic and avoids leftover/duplicated runs.

$d = $word.ActiveDocument

function New-PkgXml([string]$bodyFragment) {
    return '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:body>' + $bodyFragment + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
}

function Find-ParagraphIndex([string]$needle) {
    $count = $d.Paragraphs.Count
    for ($i = 1; $i -le $count; $i++) {
        $t = $d.Paragraphs.Item($i).Range.Text
        if ($t -like "*$needle*") {
            return $i
        }
    }
    throw "Paragraph containing '$needle' not found"
}

# ---------------------------------------------------------------------------
# 1) "... ver donde se trunca el numero."  -> add gramStart/gramEnd around "numero"
# ---------------------------------------------------------------------------
$idx1 = Find-ParagraphIndex "se trunca el numero"
$rng1 = $d.Paragraphs.Item($idx1).Range
$xml1 = '<w:p w:rsidR="004065DA" w:rsidRDefault="004065DA" w:rsidP="00B85F29">' +
    '<w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr>' +
    '<w:r><w:t xml:space="preserve">Los puntajes muestran solo </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/><w:r><w:t>int</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> y no </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/><w:r><w:t>double</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve">, ver donde se trunca el </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/>' +
    '<w:r><w:t>numero</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/>' +
    '<w:r><w:t>.</w:t></w:r>' +
    '</w:p>'
$rng1.InsertXML((New-PkgXml $xml1))

# ---------------------------------------------------------------------------
# 2) "... ver tema del token." -> split "token" off with spellStart/spellEnd and
#    insert a new paragraph "FILTRO POR ZONA A LA HORA DE BUSCAR." right after it
#    (the _GoBack bookmark now lives on the new paragraph).
# ---------------------------------------------------------------------------
$idx2 = Find-ParagraphIndex "ver tema del token"
$rng2 = $d.Paragraphs.Item($idx2).Range
$xml2 = '<w:p w:rsidR="000D22DA" w:rsidRPr="00783F18" w:rsidRDefault="000D22DA" w:rsidP="009F43A7">' +
    '<w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr>' +
    '<w:r><w:t xml:space="preserve">Ver registro de cliente, da error al ingresar cliente no existe el servicio, ver tema del </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/><w:r><w:t>token</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t>.</w:t></w:r>' +
    '</w:p>' +
    '<w:p>' +
    '<w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr>' +
    '<w:r><w:t>FILTRO POR ZONA A LA HORA DE BUSCAR.</w:t></w:r>' +
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' +
    '</w:p>'
$rng2.InsertXML((New-PkgXml $xml2))

# ---------------------------------------------------------------------------
# 3) "... Falta mostrar mas notificaciones." -> add gramStart/gramEnd around "mas"
# ---------------------------------------------------------------------------
$idx3 = Find-ParagraphIndex "Falta mostrar mas notificaciones"
$rng3 = $d.Paragraphs.Item($idx3).Range
$xml3 = '<w:p w:rsidR="00E221A6" w:rsidRDefault="00E221A6" w:rsidP="00E221A6">' +
    '<w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="8"/></w:numPr></w:pPr>' +
    '<w:r w:rsidRPr="009F43A7"><w:rPr><w:highlight w:val="yellow"/></w:rPr><w:t>Notificar y mostrar datos del due' + [char]0xF1 + 'o de la publicaci' + [char]0xF3 + 'n al trabajador contratado</w:t></w:r>' +
    '<w:r><w:t>.</w:t></w:r>' +
    '<w:r w:rsidR="009F43A7"><w:t xml:space="preserve"> Falta mostrar </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/>' +
    '<w:r w:rsidR="009F43A7"><w:t>mas</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/>' +
    '<w:r w:rsidR="009F43A7"><w:t xml:space="preserve"> notificaciones.</w:t></w:r>' +
    '</w:p>'
$rng3.InsertXML((New-PkgXml $xml3))

# ---------------------------------------------------------------------------
# 4) heading "listado-solicitudes-ofrecidas" -> wrap whole run with gramStart/gramEnd
# ---------------------------------------------------------------------------
$idx4 = Find-ParagraphIndex "listado-solicitudes-ofrecidas"
$rng4 = $d.Paragraphs.Item($idx4).Range
$xml4 = '<w:p w:rsidR="00B85F29" w:rsidRDefault="00421F4C" w:rsidP="00421F4C">' +
    '<w:pPr><w:pStyle w:val="Ttulo2"/></w:pPr>' +
    '<w:proofErr w:type="gramStart"/>' +
    '<w:r w:rsidRPr="00421F4C"><w:t>listado-solicitudes-ofrecidas</w:t></w:r>' +
    '<w:proofErr w:type="gramEnd"/>' +
    '</w:p>'
$rng4.InsertXML((New-PkgXml $xml4))

Write-Output "Done. Paragraph count: $($d.Paragraphs.Count)"
